$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename columns, add a new third header ---
$ws.Range("A1").Value = "Employee_ID"
$ws.Range("B1").Value = "Phone_Number"
$ws.Range("C1").Value = "E_mail"

# Give the new header cell the same look (bold/border/centered) as the
# existing header cells, by copying A1's format onto it.
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# --- Data row 2 ---
$ws.Range("A2").Value = 870840
$ws.Range("B2").Value = 1234567890

# --- Data row 3 ---
$ws.Range("A3").Value = 873843
# This phone number is stored as text rather than a number. Build it with
# TEXT() and then flatten the formula down to its literal string result so
# the cell ends up holding a plain text value (not a live formula).
$ws.Range("B3").Formula = "=TEXT(1234567890,""0"")"
$ws.Range("B3").Copy()
$ws.Range("B3").PasteSpecial(-4163)

# --- E-mail column, filled in after the phone numbers ---
$ws.Range("C2").Value = "gui_testing@gmail.com"
$ws.Range("C3").Value = "gui_new_login@testing.com"

$excel.CutCopyMode = 0
